# ---------------------------------------------------------------------------
# Köv.lista c. dia tartalmi elemeinek megírása
#
# 1) Slide 3 ("Jelenlegi helyzet"): merge the three trailing runs of the
#    content placeholder's paragraph into a single run.
# 2) Add a new slide 4 ("Követelménylista") with a title + bulleted content
#    placeholder.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Slide 3: merge runs ------------------------------------------------

$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

# The paragraph currently looks like:
#   "Cégünk ... Önök " + "cége" + " készítette ... folyamatainkat" + ", ezért megkérjük " + "önöket ... továbbfejlesztésére."
# Runs 3+4+5 (character offset 77, length 623 out of 699) need to collapse
# into a single run while keeping run 3's formatting.
$merged = $body3.Characters(77, 623)
$merged.Text = " készítette el számunkra. A programmal elégedettek vagyunk, azonban a gyorsan fejlődő világban, rohamosan változó piaci helyzet mellett, elkerülhetetlenné vált ezen program továbbfejlesztése is. Programunkban képesek vagyunk rögzíteni az autókat valamint a hozzájuk kapcsolódó ügyfeleket. Itt követjük nyomon a szerelések árát is. A kifizetett autókat egy hónapon belül töröljük az adatbázisból. Azt is tudni érdemes, hogy az Önök által használt alkalmazásának is vannak hiányosságai amik orvosolása nagyban megkönnyítené a jelenlegi adminisztrációs folyamatainkat, ezért megkérjük önöket az alkalmazás továbbfejlesztésére."

# --- 2) Add slide 4: "Követelménylista" ------------------------------------

$s4 = $p.Slides.Add(4, 2)

$title4 = $s4.Shapes.Item(1).TextFrame.TextRange
$title4.Text = "Követelménylista"
$title4.LanguageID = "hu-HU"

$body4 = $s4.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1
$body4.Text = "A rendszer legyen képes az autók mellett a szerelők nyilvántartására is."
$body4.LanguageID = "hu-HU"

# Paragraph 2
$body4.InsertAfter("`rA rendszer legyen képes részletes adattárolásra (egy tulajdonoshoz több autót is hozzá lehessen rendelni.)")
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).LanguageID = "hu-HU"

# Paragraph 3 (built from three runs)
$body4.InsertAfter("`rLegyünk ")
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).LanguageID = "hu-HU"
$body4.InsertAfter("képesek külön-külön ")
$body4.InsertAfter("szerkeszteni az ilyen jellegű autókat (melyek egy tulajdonoshoz vannak rendelve.)")

# Paragraph 4
$body4.InsertAfter("`rA programunk továbbfejlesztett verziója is, egyszerű, letisztult, könnyedén kezelhető, felhasználóbarát felülettel rendelkezzen.")
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).LanguageID = "hu-HU"

# Paragraph 5 (trailing empty, no bullet)
$body4.InsertAfter("`r")
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).LanguageID = "hu-HU"
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).ParagraphFormat.Bullet.Visible = 0
